{"js": "// Add alt-text (\"descr\") to the QR code image in the document.\n// The document contains a single inline picture (the QR code); set its\n// altTextDescription, which Word serializes as the `descr` attribute on\n// both wp:docPr and pic:cNvPr.\nconst pictures = context.document.body.inlinePictures;\npictures.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < pictures.items.length; i++) {\n  pictures.items[i].altTextDescription = \"Survey QR Code\";\n}\n\nawait context.sync();\n", "ps1": "# Add alt-text (\"descr\") to the QR code image in the document.\n# The document contains a single inline picture (the QR code); set its\n# AlternativeText, which Word serializes as the `descr` attribute on\n# both wp:docPr and pic:cNvPr.\n$d = $word.ActiveDocument\nforeach ($shp in $d.InlineShapes) {\n    $shp.AlternativeText = \"Survey QR Code\"\n}\n"}
